$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username values first ...
$ws.Range("C2").Value = "shankar1239"
$ws.Range("C3").Value = "prem1239"

# ... then the email addresses, matching the order new shared strings
# were appended to the workbook's string table.
$ws.Range("B2").Value = "shankar87@gmail.com"
$ws.Range("B3").Value = "prem38@gmail.com"

# Move / leave the selection where the author left it on save
$ws.Range("E17").Select()
